$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

$ws.Range("A32:G32").Copy()
$ws.Range("A52:G52").PasteSpecial(-4122)

$ws.Range("A52").Value = "Linking_AutoUser"
$ws.Range("B52").Value = "Password1"
$ws.Range("E52").Value = "Default user for Linking tests"
$ws.Range("F52").Value = "N"
$ws.Range("G52").Value = "linking.autouser@mailinator.com"
